$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target roster data (player, positions, team) in final row order
$data = @(
    @('Ja Morant', 'PG', 'Memphis Grizzlies'),
    @('Mikal Bridges', 'SG,SF,PF', 'New York Knicks'),
    @('Josh Giddey', 'PG,SG,SF', 'Chicago Bulls'),
    @('Miles Bridges', 'SF,PF', 'Charlotte Hornets'),
    @('Dillon Brooks', 'SG,SF', 'Houston Rockets'),
    @('Kevon Looney', 'PF,C', 'Golden State Warriors'),
    @('Nikola Vucevic', 'PF,C', 'Chicago Bulls'),
    @('Nick Richards', 'C', 'Phoenix Suns'),
    @('Brook Lopez', 'C', 'Milwaukee Bucks'),
    @('Aaron Gordon', 'PF,C', 'Denver Nuggets'),
    @('DeMar DeRozan', 'SF,PF', 'Sacramento Kings'),
    @('Tyler Herro', 'PG,SG', 'Miami Heat'),
    @('De''Aaron Fox', 'PG', 'Sacramento Kings'),
    @('Shaedon Sharpe', 'SG,SF', 'Portland Trail Blazers'),
    @('Scottie Barnes', 'PG,SG,SF,PF', 'Toronto Raptors'),
    @('Luka Doncic', 'PG,SG', 'Dallas Mavericks'),
    @('Evan Mobley', 'PF,C', 'Cleveland Cavaliers'),
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}

